$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 253.875
$ws.Cells.Item(2, 9).Value = 128.2
$ws.Cells.Item(2, 11).Value = 128.2
$ws.Cells.Item(2, 13).Value = -15.19999999999999
$ws.Cells.Item(6, 8).Value = 1320.2632
$ws.Cells.Item(6, 9).Value = 249.66667
$ws.Cells.Item(6, 10).Value = 1814.3846
$ws.Cells.Item(6, 11).Value = 749.00001
$ws.Cells.Item(6, 12).Value = 5443.1538
$ws.Cells.Item(6, 13).Value = -637.00001
$ws.Cells.Item(6, 14).Value = -5667.1538
$ws.Cells.Item(12, 8).Value = 516.8
$ws.Cells.Item(12, 9).Value = 239.875
$ws.Cells.Item(12, 10).Value = 1624.5
$ws.Cells.Item(12, 11).Value = 239.875
$ws.Cells.Item(12, 12).Value = 1624.5
$ws.Cells.Item(12, 13).Value = -69.875
$ws.Cells.Item(12, 14).Value = -1964.5
$ws.Cells.Item(15, 8).Value = 462515.25
$ws.Cells.Item(15, 9).Value = 462515.25
$ws.Cells.Item(15, 11).Value = 1387545.75
$ws.Cells.Item(15, 13).Value = -1387376.75
$ws.Cells.Item(32, 8).Value = 1180.2
$ws.Cells.Item(32, 10).Value = 1250.25
$ws.Cells.Item(32, 12).Value = 1250.25
$ws.Cells.Item(32, 14).Value = -1902.25
$ws.Cells.Item(40, 8).Value = 3994.125
$ws.Cells.Item(40, 9).Value = 2809.6365
$ws.Cells.Item(40, 11).Value = 2809.6365
$ws.Cells.Item(40, 13).Value = -2634.6365
$ws.Cells.Item(42, 8).Value = 221.15384
$ws.Cells.Item(42, 9).Value = 86
$ws.Cells.Item(42, 11).Value = 258
$ws.Cells.Item(42, 13).Value = -28
$ws.Cells.Item(43, 8).Value = 10235.615
$ws.Cells.Item(43, 10).Value = 12296.111
$ws.Cells.Item(43, 12).Value = 12296.111
$ws.Cells.Item(43, 14).Value = -12434.111
$ws.Cells.Item(51, 8).Value = 38231.645
$ws.Cells.Item(51, 9).Value = 2535.25
$ws.Cells.Item(51, 10).Value = 52510.2
$ws.Cells.Item(51, 11).Value = 2535.25
$ws.Cells.Item(51, 12).Value = 52510.2
$ws.Cells.Item(51, 13).Value = -2051.25
$ws.Cells.Item(51, 14).Value = -53478.2
$ws.Cells.Item(53, 8).Value = 455.14285
$ws.Cells.Item(53, 9).Value = 329.25
$ws.Cells.Item(53, 10).Value = 623
$ws.Cells.Item(53, 11).Value = 329.25
$ws.Cells.Item(53, 12).Value = 623
$ws.Cells.Item(53, 13).Value = 307.75
$ws.Cells.Item(53, 14).Value = -1897
$ws.Cells.Item(62, 8).Value = 76931816
$ws.Cells.Item(62, 9).Value = 111119240
$ws.Cells.Item(62, 11).Value = 111119240
$ws.Cells.Item(62, 13).Value = -111118616
$ws.Cells.Item(65, 8).Value = 76931816
$ws.Cells.Item(65, 9).Value = 111119240
$ws.Cells.Item(65, 11).Value = 555596200
$ws.Cells.Item(65, 13).Value = -555593080
$ws.Cells.Item(86, 8).Value = 58856660
$ws.Cells.Item(86, 9).Value = 3631.5386
$ws.Cells.Item(86, 10).Value = 250129000
$ws.Cells.Item(86, 11).Value = 3631.5386
$ws.Cells.Item(86, 12).Value = 250129000
$ws.Cells.Item(86, 13).Value = -2508.5386
$ws.Cells.Item(86, 14).Value = -250131246
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).Value = $null
$ws.Cells.Item(89, 8).Value = 58856660
$ws.Cells.Item(89, 9).Value = 3631.5386
$ws.Cells.Item(89, 10).Value = 250129000
$ws.Cells.Item(89, 11).Value = 18157.693
$ws.Cells.Item(89, 12).Value = 1250645000
$ws.Cells.Item(89, 13).Value = -12541.693
$ws.Cells.Item(89, 14).Value = -1250656232
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).Value = $null
$ws.Cells.Item(99, 8).Value = 1823605.1
$ws.Cells.Item(99, 10).Value = 2055.2
$ws.Cells.Item(99, 12).Value = 6165.599999999999
$ws.Cells.Item(99, 14).Value = -9161.599999999999
$ws.Cells.Item(113, 8).Value = 18079.818
$ws.Cells.Item(113, 9).Value = 18859.875
$ws.Cells.Item(113, 11).Value = 18859.875
$ws.Cells.Item(113, 13).Value = -15605.875
$ws.Cells.Item(125, 8).Value = 9369.799999999999
$ws.Cells.Item(125, 9).Value = 32000
$ws.Cells.Item(125, 10).Value = 3712.25
$ws.Cells.Item(125, 11).Value = 288000
$ws.Cells.Item(125, 12).Value = 33410.25
$ws.Cells.Item(125, 13).Value = -285540
$ws.Cells.Item(125, 14).Value = -38330.25
$ws.Cells.Item(129, 8).Value = 2552.6
$ws.Cells.Item(129, 9).Value = 1682.7142
$ws.Cells.Item(129, 10).Value = 4582.3335
$ws.Cells.Item(129, 11).Value = 5048.142599999999
$ws.Cells.Item(129, 12).Value = 13747.0005
$ws.Cells.Item(129, 13).Value = -48.14259999999922
$ws.Cells.Item(129, 14).Value = -23747.0005
$ws.Cells.Item(135, 8).Value = 3828.85
$ws.Cells.Item(135, 9).Value = 3916.9412
$ws.Cells.Item(135, 10).Value = 3329.6667
$ws.Cells.Item(135, 11).Value = 35252.4708
$ws.Cells.Item(135, 12).Value = 29967.0003
$ws.Cells.Item(135, 13).Value = -32717.4708
$ws.Cells.Item(135, 14).Value = -35037.0003
$ws.Cells.Item(137, 8).Value = 3581.625
$ws.Cells.Item(137, 9).Value = 3051.889
$ws.Cells.Item(137, 10).Value = 4262.7144
$ws.Cells.Item(137, 11).Value = 9155.667000000001
$ws.Cells.Item(137, 12).Value = 12788.1432
$ws.Cells.Item(137, 13).Value = -6605.667000000001
$ws.Cells.Item(137, 14).Value = -17888.1432
$ws.Cells.Item(138, 8).Value = 142923.72
$ws.Cells.Item(138, 9).Value = 1465641.5
$ws.Cells.Item(138, 10).Value = 6090.8506
$ws.Cells.Item(138, 11).Value = 4396924.5
$ws.Cells.Item(138, 12).Value = 18272.5518
$ws.Cells.Item(138, 13).Value = -4391784.5
$ws.Cells.Item(138, 14).Value = -28552.5518
$ws.Cells.Item(139, 8).Value = 167301.25
$ws.Cells.Item(139, 10).Value = 167301.25
$ws.Cells.Item(139, 12).Value = 167301.25
$ws.Cells.Item(139, 14).Value = -177581.25
$ws.Cells.Item(140, 8).Value = 99999
$ws.Cells.Item(140, 10).Value = 99999
$ws.Cells.Item(140, 12).Value = 99999
$ws.Cells.Item(140, 14).Value = -110359
$ws.Cells.Item(141, 8).Value = 6527.9414
$ws.Cells.Item(141, 9).Value = 6685.9375
$ws.Cells.Item(141, 11).Value = 20057.8125
$ws.Cells.Item(141, 13).Value = -14877.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 7110.75
$ws.Cells.Item(3, 9).Value = 6147.6665
$ws.Cells.Item(3, 11).Value = 6147.6665
$ws.Cells.Item(3, 13).Value = -6032.6665
$ws.Cells.Item(32, 8).Value = 14871.759
$ws.Cells.Item(32, 9).Value = 13434.66
$ws.Cells.Item(32, 10).Value = 32835.5
$ws.Cells.Item(32, 11).Value = 13434.66
$ws.Cells.Item(32, 12).Value = 32835.5
$ws.Cells.Item(32, 13).Value = -13147.66
$ws.Cells.Item(32, 14).Value = -33409.5
$ws.Cells.Item(41, 8).Value = 3920.6667
$ws.Cells.Item(41, 9).Value = 2508
$ws.Cells.Item(41, 10).Value = 5333.3335
$ws.Cells.Item(41, 11).Value = 2508
$ws.Cells.Item(41, 12).Value = 5333.3335
$ws.Cells.Item(41, 13).Value = -2094
$ws.Cells.Item(41, 14).Value = -6161.3335
$ws.Cells.Item(45, 8).Value = 4235.4194
$ws.Cells.Item(45, 9).Value = 4154.1577
$ws.Cells.Item(45, 11).Value = 4154.1577
$ws.Cells.Item(45, 13).Value = -3777.1577
$ws.Cells.Item(61, 8).Value = 8971.566000000001
$ws.Cells.Item(61, 9).Value = 9412.888999999999
$ws.Cells.Item(61, 11).Value = 9412.888999999999
$ws.Cells.Item(61, 13).Value = -9200.888999999999
$ws.Cells.Item(74, 8).Value = 4428.8667
$ws.Cells.Item(74, 9).Value = 60649
$ws.Cells.Item(74, 10).Value = 1813.9767
$ws.Cells.Item(74, 11).Value = 60649
$ws.Cells.Item(74, 12).Value = 1813.9767
$ws.Cells.Item(74, 13).Value = -59775
$ws.Cells.Item(74, 14).Value = -3561.9767
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).Value = $null
$ws.Cells.Item(77, 8).Value = 4428.8667
$ws.Cells.Item(77, 9).Value = 60649
$ws.Cells.Item(77, 10).Value = 1813.9767
$ws.Cells.Item(77, 11).Value = 303245
$ws.Cells.Item(77, 12).Value = 9069.8835
$ws.Cells.Item(77, 13).Value = -298877
$ws.Cells.Item(77, 14).Value = -17805.8835
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).Value = $null
$ws.Cells.Item(82, 8).Value = 44164
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).Value = $null
$ws.Cells.Item(85, 8).Value = 44164
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).Value = $null
$ws.Cells.Item(102, 8).Value = 8777.121999999999
$ws.Cells.Item(102, 9).Value = 9624.593000000001
$ws.Cells.Item(102, 10).Value = 7142.7144
$ws.Cells.Item(102, 11).Value = 9624.593000000001
$ws.Cells.Item(102, 12).Value = 7142.7144
$ws.Cells.Item(102, 13).Value = -8002.593000000001
$ws.Cells.Item(102, 14).Value = -10386.7144
$ws.Cells.Item(122, 8).Value = 1204594.6
$ws.Cells.Item(122, 9).Value = 4611.4165
$ws.Cells.Item(122, 11).Value = 13834.2495
$ws.Cells.Item(122, 13).Value = -11384.2495
$ws.Cells.Item(125, 8).Value = 91079.336
$ws.Cells.Item(125, 10).Value = 91079.336
$ws.Cells.Item(125, 12).Value = 91079.336
$ws.Cells.Item(125, 14).Value = -100919.336
$ws.Cells.Item(132, 8).Value = 1914.0625
$ws.Cells.Item(132, 9).Value = 1366.88
$ws.Cells.Item(132, 11).Value = 4100.64
$ws.Cells.Item(132, 13).Value = -1570.64
$ws.Cells.Item(136, 8).Value = 8971.566000000001
$ws.Cells.Item(136, 9).Value = 9412.888999999999
$ws.Cells.Item(136, 11).Value = 28238.667
$ws.Cells.Item(136, 13).Value = -25688.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 10997.5
$ws.Cells.Item(8, 9).Value = 7000
$ws.Cells.Item(8, 11).Value = 7000
$ws.Cells.Item(8, 13).Value = -6860
$ws.Cells.Item(20, 8).Value = 5126
$ws.Cells.Item(20, 9).Value = 3871.1428
$ws.Cells.Item(20, 10).Value = 7322
$ws.Cells.Item(20, 11).Value = 3871.1428
$ws.Cells.Item(20, 12).Value = 7322
$ws.Cells.Item(20, 13).Value = -3624.1428
$ws.Cells.Item(20, 14).Value = -7816
$ws.Cells.Item(26, 8).Value = 125999.6
$ws.Cells.Item(26, 10).Value = 109999
$ws.Cells.Item(26, 12).Value = 109999
$ws.Cells.Item(26, 14).Value = -110583
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 13).Value = $null
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).Value = $null
$ws.Cells.Item(61, 8).Value = 149998
$ws.Cells.Item(61, 10).Value = 149998
$ws.Cells.Item(61, 12).Value = 149998
$ws.Cells.Item(61, 14).Value = -150624
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).Value = $null
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).Value = $null
$ws.Cells.Item(80, 8).Value = 442.88235
$ws.Cells.Item(80, 10).Value = 391.69232
$ws.Cells.Item(80, 12).Value = 391.69232
$ws.Cells.Item(80, 14).Value = -2387.69232
$ws.Cells.Item(82, 8).Value = 9876.666999999999
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).Value = $null
$ws.Cells.Item(83, 8).Value = 442.88235
$ws.Cells.Item(83, 10).Value = 391.69232
$ws.Cells.Item(83, 12).Value = 1958.4616
$ws.Cells.Item(83, 14).Value = -11942.4616
$ws.Cells.Item(85, 8).Value = 9876.666999999999
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).Value = $null
$ws.Cells.Item(94, 8).Value = 2727.6572
$ws.Cells.Item(94, 9).Value = 1985.7
$ws.Cells.Item(94, 10).Value = 7179.4
$ws.Cells.Item(94, 11).Value = 1985.7
$ws.Cells.Item(94, 12).Value = 7179.4
$ws.Cells.Item(94, 13).Value = -1534.7
$ws.Cells.Item(94, 14).Value = -8081.4
$ws.Cells.Item(105, 8).Value = 9029.35
$ws.Cells.Item(105, 9).Value = 9184.786
$ws.Cells.Item(105, 11).Value = 9184.786
$ws.Cells.Item(105, 13).Value = -7437.786
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 14).Value = $null
$ws.Cells.Item(123, 8).Value = 96992
$ws.Cells.Item(123, 10).Value = 96992
$ws.Cells.Item(123, 12).Value = 96992
$ws.Cells.Item(123, 14).Value = -106792
$ws.Cells.Item(133, 8).Value = 75000
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 75000
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 75000
$ws.Cells.Item(133, 13).Value = $null
$ws.Cells.Item(133, 14).Value = -85120
$ws.Cells.Item(134, 8).Value = 4343.077
$ws.Cells.Item(134, 9).Value = 3432.5
$ws.Cells.Item(134, 10).Value = 5800
$ws.Cells.Item(134, 11).Value = 10297.5
$ws.Cells.Item(134, 12).Value = 17400
$ws.Cells.Item(134, 13).Value = -7762.5
$ws.Cells.Item(134, 14).Value = -22470

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4131
$ws.Cells.Item(31, 9).Value = 2403.6428
$ws.Cells.Item(31, 10).Value = 5138.625
$ws.Cells.Item(31, 11).Value = 2403.6428
$ws.Cells.Item(31, 12).Value = 5138.625
$ws.Cells.Item(31, 13).Value = -2108.6428
$ws.Cells.Item(31, 14).Value = -5728.625
$ws.Cells.Item(34, 8).Value = 4131
$ws.Cells.Item(34, 9).Value = 2403.6428
$ws.Cells.Item(34, 10).Value = 5138.625
$ws.Cells.Item(34, 11).Value = 2403.6428
$ws.Cells.Item(34, 12).Value = 5138.625
$ws.Cells.Item(34, 13).Value = -2201.6428
$ws.Cells.Item(34, 14).Value = -5542.625
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 13).Value = $null
$ws.Cells.Item(54, 14).Value = $null
$ws.Cells.Item(58, 8).Value = 3189.52
$ws.Cells.Item(58, 9).Value = 2612.6875
$ws.Cells.Item(58, 11).Value = 2612.6875
$ws.Cells.Item(58, 13).Value = -2409.6875
$ws.Cells.Item(68, 8).Value = 10000
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = $null
$ws.Cells.Item(71, 8).Value = 10000
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = $null
$ws.Cells.Item(74, 8).Value = 25000
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).Value = $null
$ws.Cells.Item(77, 8).Value = 25000
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).Value = $null
$ws.Cells.Item(97, 8).Value = 39033
$ws.Cells.Item(97, 10).Value = 39033
$ws.Cells.Item(97, 12).Value = 39033
$ws.Cells.Item(97, 14).Value = -41015
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 13).Value = $null
$ws.Cells.Item(108, 14).Value = $null
$ws.Cells.Item(109, 8).Value = 149998
$ws.Cells.Item(109, 10).Value = 149998
$ws.Cells.Item(109, 12).Value = 149998
$ws.Cells.Item(109, 14).Value = -152078
$ws.Cells.Item(114, 8).Value = 55553.332
$ws.Cells.Item(114, 10).Value = 55553.332
$ws.Cells.Item(114, 12).Value = 55553.332
$ws.Cells.Item(114, 14).Value = -64231.332
$ws.Cells.Item(122, 8).Value = 35000
$ws.Cells.Item(122, 9).Value = 62500
$ws.Cells.Item(122, 11).Value = 187500
$ws.Cells.Item(122, 13).Value = -185050
$ws.Cells.Item(132, 8).Value = 19980.818
$ws.Cells.Item(132, 9).Value = 30141.285
$ws.Cells.Item(132, 11).Value = 90423.855
$ws.Cells.Item(132, 13).Value = -87893.855
$ws.Cells.Item(134, 8).Value = 1877.7307
$ws.Cells.Item(134, 9).Value = 1867
$ws.Cells.Item(134, 11).Value = 5601
$ws.Cells.Item(134, 13).Value = -3066
$ws.Cells.Item(136, 8).Value = 3189.52
$ws.Cells.Item(136, 9).Value = 2612.6875
$ws.Cells.Item(136, 11).Value = 7838.0625
$ws.Cells.Item(136, 13).Value = -5288.0625
$ws.Cells.Item(141, 8).Value = 1120099
$ws.Cells.Item(141, 10).Value = 1120099
$ws.Cells.Item(141, 12).Value = 1120099
$ws.Cells.Item(141, 14).Value = -1130459

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 92.59999999999999
$ws.Cells.Item(2, 9).Value = 17
$ws.Cells.Item(2, 10).Value = 243.8
$ws.Cells.Item(2, 11).Value = 102
$ws.Cells.Item(2, 12).Value = 1462.8
$ws.Cells.Item(2, 13).Value = 11
$ws.Cells.Item(2, 14).Value = -1688.8
$ws.Cells.Item(4, 8).Value = 43306064
$ws.Cells.Item(4, 9).Value = 39387036
$ws.Cells.Item(4, 10).Value = 60101890
$ws.Cells.Item(4, 11).Value = 118161108
$ws.Cells.Item(4, 12).Value = 180305670
$ws.Cells.Item(4, 13).Value = -118160996
$ws.Cells.Item(4, 14).Value = -180305894
$ws.Cells.Item(12, 8).Value = 55.923077
$ws.Cells.Item(12, 9).Value = 90
$ws.Cells.Item(12, 11).Value = 270
$ws.Cells.Item(12, 13).Value = -97
$ws.Cells.Item(14, 8).Value = 578.8
$ws.Cells.Item(14, 9).Value = 578.8
$ws.Cells.Item(14, 11).Value = 1736.4
$ws.Cells.Item(14, 13).Value = -1563.4
$ws.Cells.Item(17, 8).Value = 2702.8
$ws.Cells.Item(17, 10).Value = 5995
$ws.Cells.Item(17, 12).Value = 17985
$ws.Cells.Item(17, 14).Value = -18323
$ws.Cells.Item(34, 8).Value = 2303.2856
$ws.Cells.Item(34, 9).Value = 2064.6
$ws.Cells.Item(34, 10).Value = 2900
$ws.Cells.Item(34, 11).Value = 6193.799999999999
$ws.Cells.Item(34, 12).Value = 8700
$ws.Cells.Item(34, 13).Value = -6109.799999999999
$ws.Cells.Item(34, 14).Value = -8868
$ws.Cells.Item(55, 8).Value = 10521.091
$ws.Cells.Item(55, 10).Value = 12761.223
$ws.Cells.Item(55, 12).Value = 38283.669
$ws.Cells.Item(55, 14).Value = -38637.669
$ws.Cells.Item(92, 8).Value = 1515.8334
$ws.Cells.Item(92, 9).Value = 1724.75
$ws.Cells.Item(92, 10).Value = 1098
$ws.Cells.Item(92, 11).Value = 5174.25
$ws.Cells.Item(92, 12).Value = 3294
$ws.Cells.Item(92, 13).Value = -3926.25
$ws.Cells.Item(92, 14).Value = -5790
$ws.Cells.Item(121, 8).Value = 515.8461
$ws.Cells.Item(121, 9).Value = 421.5
$ws.Cells.Item(121, 10).Value = 830.3333
$ws.Cells.Item(121, 11).Value = 1264.5
$ws.Cells.Item(121, 12).Value = 2490.9999
$ws.Cells.Item(121, 13).Value = 45.5
$ws.Cells.Item(121, 14).Value = -5110.9999
$ws.Cells.Item(122, 8).Value = 5817.8066
$ws.Cells.Item(122, 10).Value = 8143.8
$ws.Cells.Item(122, 12).Value = 73294.2
$ws.Cells.Item(122, 14).Value = -78194.2
$ws.Cells.Item(132, 8).Value = 25301.26
$ws.Cells.Item(132, 9).Value = 1233.1666
$ws.Cells.Item(132, 10).Value = 32177.857
$ws.Cells.Item(132, 11).Value = 11098.4994
$ws.Cells.Item(132, 12).Value = 289600.713
$ws.Cells.Item(132, 13).Value = -8568.499400000001
$ws.Cells.Item(132, 14).Value = -294660.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 78.08696
$ws.Cells.Item(2, 9).Value = 66.47619
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(2, 11).Value = 66.47619
$ws.Cells.Item(2, 12).Value = 200
$ws.Cells.Item(2, 13).Value = 46.52381
$ws.Cells.Item(2, 14).Value = -426
$ws.Cells.Item(40, 8).Value = 9997.25
$ws.Cells.Item(40, 9).Value = 9994
$ws.Cells.Item(40, 10).Value = 9998.333000000001
$ws.Cells.Item(40, 11).Value = 9994
$ws.Cells.Item(40, 12).Value = 9998.333000000001
$ws.Cells.Item(40, 13).Value = -9843
$ws.Cells.Item(40, 14).Value = -10300.333
$ws.Cells.Item(52, 8).Value = 24874.75
$ws.Cells.Item(52, 9).Value = 35000
$ws.Cells.Item(52, 10).Value = 21499.666
$ws.Cells.Item(52, 11).Value = 35000
$ws.Cells.Item(52, 12).Value = 21499.666
$ws.Cells.Item(52, 13).Value = -34741
$ws.Cells.Item(52, 14).Value = -22017.666
$ws.Cells.Item(55, 8).Value = 10484
$ws.Cells.Item(55, 10).Value = 12397.5
$ws.Cells.Item(55, 12).Value = 12397.5
$ws.Cells.Item(55, 14).Value = -13051.5
$ws.Cells.Item(70, 8).Value = 6647.2104
$ws.Cells.Item(70, 9).Value = 5835.143
$ws.Cells.Item(70, 11).Value = 5835.143
$ws.Cells.Item(70, 13).Value = -5565.143
$ws.Cells.Item(73, 8).Value = 6647.2104
$ws.Cells.Item(73, 9).Value = 5835.143
$ws.Cells.Item(73, 11).Value = 5835.143
$ws.Cells.Item(73, 13).Value = -4899.143
$ws.Cells.Item(80, 8).Value = 17288.2
$ws.Cells.Item(80, 9).Value = 23191.715
$ws.Cells.Item(80, 11).Value = 23191.715
$ws.Cells.Item(80, 13).Value = -22193.715
$ws.Cells.Item(83, 8).Value = 17288.2
$ws.Cells.Item(83, 9).Value = 23191.715
$ws.Cells.Item(83, 11).Value = 115958.575
$ws.Cells.Item(83, 13).Value = -110966.575
$ws.Cells.Item(93, 8).Value = 33888.5
$ws.Cells.Item(93, 10).Value = 33888.5
$ws.Cells.Item(93, 12).Value = 33888.5
$ws.Cells.Item(93, 14).Value = -37632.5
$ws.Cells.Item(97, 8).Value = 11449
$ws.Cells.Item(97, 9).Value = 13076.667
$ws.Cells.Item(97, 11).Value = 13076.667
$ws.Cells.Item(97, 13).Value = -12580.667
$ws.Cells.Item(113, 8).Value = 19430.143
$ws.Cells.Item(113, 9).Value = 39337
$ws.Cells.Item(113, 10).Value = 4500
$ws.Cells.Item(113, 11).Value = 39337
$ws.Cells.Item(113, 12).Value = 4500
$ws.Cells.Item(113, 13).Value = -37167
$ws.Cells.Item(113, 14).Value = -8840
$ws.Cells.Item(132, 8).Value = 4919.5625
$ws.Cells.Item(132, 9).Value = 5001.423
$ws.Cells.Item(132, 10).Value = 4564.8335
$ws.Cells.Item(132, 11).Value = 15004.269
$ws.Cells.Item(132, 12).Value = 13694.5005
$ws.Cells.Item(132, 13).Value = -12474.269
$ws.Cells.Item(132, 14).Value = -18754.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 720
$ws.Cells.Item(22, 9).Value = 692.3077
$ws.Cells.Item(22, 10).Value = 900
$ws.Cells.Item(22, 11).Value = 692.3077
$ws.Cells.Item(22, 12).Value = 900
$ws.Cells.Item(22, 13).Value = -397.3077
$ws.Cells.Item(22, 14).Value = -1490
$ws.Cells.Item(27, 8).Value = 720
$ws.Cells.Item(27, 9).Value = 692.3077
$ws.Cells.Item(27, 10).Value = 900
$ws.Cells.Item(27, 11).Value = 692.3077
$ws.Cells.Item(27, 12).Value = 900
$ws.Cells.Item(27, 13).Value = -585.3077
$ws.Cells.Item(27, 14).Value = -1114
$ws.Cells.Item(46, 8).Value = 1413.2273
$ws.Cells.Item(46, 10).Value = 1757.6666
$ws.Cells.Item(46, 12).Value = 1757.6666
$ws.Cells.Item(46, 14).Value = -2133.6666
$ws.Cells.Item(55, 8).Value = 3408.8
$ws.Cells.Item(55, 9).Value = 924.5
$ws.Cells.Item(55, 10).Value = 5065
$ws.Cells.Item(55, 11).Value = 924.5
$ws.Cells.Item(55, 12).Value = 5065
$ws.Cells.Item(55, 13).Value = -751.5
$ws.Cells.Item(55, 14).Value = -5411
$ws.Cells.Item(63, 8).Value = 54998.5
$ws.Cells.Item(63, 10).Value = 54998.5
$ws.Cells.Item(63, 12).Value = 54998.5
$ws.Cells.Item(63, 14).Value = -56496.5
$ws.Cells.Item(66, 8).Value = 54998.5
$ws.Cells.Item(66, 10).Value = 54998.5
$ws.Cells.Item(66, 12).Value = 164995.5
$ws.Cells.Item(66, 14).Value = -172483.5
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).Value = $null
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).Value = $null
$ws.Cells.Item(93, 8).Value = 4392.6
$ws.Cells.Item(93, 9).Value = 4642.0713
$ws.Cells.Item(93, 11).Value = 4642.0713
$ws.Cells.Item(93, 13).Value = -3394.0713
$ws.Cells.Item(109, 8).Value = 63323.332
$ws.Cells.Item(109, 10).Value = 63323.332
$ws.Cells.Item(109, 12).Value = 63323.332
$ws.Cells.Item(109, 14).Value = -66097.33199999999
$ws.Cells.Item(122, 8).Value = 6641.4546
$ws.Cells.Item(122, 9).Value = 6824.5386
$ws.Cells.Item(122, 10).Value = 6377
$ws.Cells.Item(122, 11).Value = 20473.6158
$ws.Cells.Item(122, 12).Value = 19131
$ws.Cells.Item(122, 13).Value = -18023.6158
$ws.Cells.Item(122, 14).Value = -24031
$ws.Cells.Item(123, 8).Value = 139984.5
$ws.Cells.Item(123, 10).Value = 139984.5
$ws.Cells.Item(123, 12).Value = 139984.5
$ws.Cells.Item(123, 14).Value = -149784.5
$ws.Cells.Item(132, 8).Value = 377449.38
$ws.Cells.Item(132, 9).Value = 601212.3
$ws.Cells.Item(132, 10).Value = 4511.1333
$ws.Cells.Item(132, 11).Value = 1803636.9
$ws.Cells.Item(132, 12).Value = 13533.3999
$ws.Cells.Item(132, 13).Value = -1801106.9
$ws.Cells.Item(132, 14).Value = -18593.3999
$ws.Cells.Item(140, 8).Value = 175666.67
$ws.Cells.Item(140, 10).Value = 175666.67
$ws.Cells.Item(140, 12).Value = 175666.67
$ws.Cells.Item(140, 14).Value = -186026.67
$ws.Cells.Item(141, 8).Value = 120178.2
$ws.Cells.Item(141, 10).Value = 120178.2
$ws.Cells.Item(141, 12).Value = 120178.2
$ws.Cells.Item(141, 14).Value = -130538.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(55, 8).Value = 5250
$ws.Cells.Item(55, 10).Value = 5250
$ws.Cells.Item(55, 12).Value = 5250
$ws.Cells.Item(55, 14).Value = -5804
$ws.Cells.Item(75, 8).Value = 70000
$ws.Cells.Item(75, 10).Value = 70000
$ws.Cells.Item(75, 12).Value = 70000
$ws.Cells.Item(75, 14).Value = -71872
$ws.Cells.Item(78, 8).Value = 70000
$ws.Cells.Item(78, 10).Value = 70000
$ws.Cells.Item(78, 12).Value = 210000
$ws.Cells.Item(78, 14).Value = -219360
$ws.Cells.Item(81, 8).Value = 22849.8
$ws.Cells.Item(81, 9).Value = 35416.332
$ws.Cells.Item(81, 11).Value = 70832.664
$ws.Cells.Item(81, 13).Value = -69771.664
$ws.Cells.Item(84, 8).Value = 22849.8
$ws.Cells.Item(84, 9).Value = 35416.332
$ws.Cells.Item(84, 11).Value = 354163.32
$ws.Cells.Item(84, 13).Value = -348859.32
$ws.Cells.Item(96, 8).Value = 5265937.5
$ws.Cells.Item(96, 9).Value = 6252716.5
$ws.Cells.Item(96, 10).Value = 3116.6667
$ws.Cells.Item(96, 11).Value = 6252716.5
$ws.Cells.Item(96, 12).Value = 3116.6667
$ws.Cells.Item(96, 13).Value = -6251343.5
$ws.Cells.Item(96, 14).Value = -5862.6667
$ws.Cells.Item(118, 8).Value = 93000
$ws.Cells.Item(118, 10).Value = 93000
$ws.Cells.Item(118, 12).Value = 93000
$ws.Cells.Item(118, 14).Value = -96314
$ws.Cells.Item(122, 8).Value = 32247.65
$ws.Cells.Item(122, 9).Value = 7045.3
$ws.Cells.Item(122, 10).Value = 57450
$ws.Cells.Item(122, 11).Value = 21135.9
$ws.Cells.Item(122, 12).Value = 172350
$ws.Cells.Item(122, 13).Value = -18685.9
$ws.Cells.Item(122, 14).Value = -177250
$ws.Cells.Item(132, 8).Value = 16410.123
$ws.Cells.Item(132, 9).Value = 18578.547
$ws.Cells.Item(132, 10).Value = 3399.5715
$ws.Cells.Item(132, 11).Value = 55735.641
$ws.Cells.Item(132, 12).Value = 10198.7145
$ws.Cells.Item(132, 13).Value = -53205.641
$ws.Cells.Item(132, 14).Value = -15258.7145
$ws.Cells.Item(136, 8).Value = 456383.4
$ws.Cells.Item(136, 9).Value = 534315.0600000001
$ws.Cells.Item(136, 10).Value = 4380
$ws.Cells.Item(136, 11).Value = 1602945.18
$ws.Cells.Item(136, 12).Value = 13140
$ws.Cells.Item(136, 13).Value = -1600395.18
$ws.Cells.Item(136, 14).Value = -18240

Write-Output "Applied 632 edits across 8 sheets"